$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "button_testResultDetails_class"
$ws.Range("B1").Value = "button_testResultDetails_class_1"
$ws.Range("C1").Value = "button_testResultDetails_class_2"
$ws.Range("D1").Value = "button_testResultDetails_internalRoleButtonName"
$ws.Range("E1").Value = "button_testResultDetails_internalRoleButtonName_1"
$ws.Range("F1").Value = "button_testResultDetails_internalRoleButtonName_2"
